$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF column (58) holds the game "Date" for each of the 30 teams (rows 2-31).
# The values were stored using the wrong format (e.g. "5-29-2011-12" instead
# of the real calendar date "2012-05-29") because of how the NBA stats site
# showed the date - fix the training data so it reflects the correct date.
$rng = $ws.Range("BF2:BF31")
$rng.NumberFormat = "@"    # keep the value literal text, not an Excel date serial

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # BF column = 58
    if ($cell.Value2 -eq "5-29-2011-12") {
        $cell.Value = "2012-05-29"
    }
}
